$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@(2000, "walkingToRunning", -61.64058685302734, 21.88947105407715, 0.6875860691070557, -3.719183167544286, 7.010346993533103, 0.6059804959730433)
    ,@(2100, "walkingToRunning", -27.69689178466797, -7.563150405883789, -14.80709171295166, -3.724935540285977, 10.48463944521817, -5.378288039294167)
    ,@(2200, "walkingToRunning", -1.68331241607666, -72.37289428710938, 70.38520812988281, 2.686929789456455, 0.002814553000709097, -3.492711760781028)
    ,@(2300, "walkingToRunning", -11.52985954284668, 10.27250671386719, 4.496110916137695, 3.455663386258163, -6.499211259321735, -1.833242598446949)
    ,@(2400, "walkingToRunning", 4.430462837219238, -20.51617813110352, 37.68183135986328, 0.03122558593747993, -2.420187681913326, 0.5764146804809926)
    ,@(2500, "walkingToRunning", -80.67318725585938, -3.729990005493164, -17.92157745361328, -7.196380597894843, -4.824799558791248, 5.860242297432641)
    ,@(2600, "walkingToRunning", -7.960968494415283, 2.627067565917969, -1.957346677780152, -1.569692446968803, -11.00871762362393, 5.813380349766132)
    ,@(2700, "walkingToRunning", -12.41016101837158, -18.66624450683594, 3.092851161956787, 3.11778094551775, -8.70344656163989, 3.647271784869103)
    ,@(2800, "walkingToRunning", 6.215286254882812, 5.965664863586426, 14.03029632568359, 4.020912179079962, 1.777516035600222, -0.3423178889535121)
    ,@(2900, "walkingToRunning", -6.658430576324463, 14.9360237121582, 14.47612380981445, 2.335198922590785, 4.166226341507659, 0.4891112804412763)
)

$startRow = 22
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
}
